$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2503.4129649087135
$ws.Range("B1").Value = 1836.6174073036473
$ws.Range("C1").Value = 1818.050812614732
$ws.Range("A2").Value = 2283.7034578987145
$ws.Range("B2").Value = 1500.5236040859768
$ws.Range("C2").Value = 1713.3126207235091
$ws.Range("A3").Value = 2560.152091236956
$ws.Range("B3").Value = 1800.3479258919385
$ws.Range("C3").Value = 1920.5342338204127
$ws.Range("A4").Value = 2537.9419460845415
$ws.Range("B4").Value = 1799.0318743020503
$ws.Range("C4").Value = 1638.4285917018726
$ws.Range("A5").Value = 2681.6693380571514
$ws.Range("B5").Value = 1744.9401510683956
$ws.Range("C5").Value = 1707.7516198242702
$ws.Range("A6").Value = 2528.0652230370224
$ws.Range("B6").Value = 1707.4757882748322
$ws.Range("C6").Value = 2003.1236656251572
$ws.Range("A7").Value = 2396.4145100733499
$ws.Range("B7").Value = 2114.840887107051
$ws.Range("C7").Value = 2076.4145363599118
$ws.Range("A8").Value = 2466.5665169556682
$ws.Range("B8").Value = 2074.1996585625466
$ws.Range("C8").Value = 2070.9993408703685
$ws.Range("A9").Value = 2646.4933228697419
$ws.Range("B9").Value = 1942.9919409553129
$ws.Range("C9").Value = 1637.3576506534309
$ws.Range("A10").Value = 2445.7574951344163
$ws.Range("B10").Value = 1437.5708577075886
$ws.Range("C10").Value = 1584.7220636104068
$ws.Range("A11").Value = 2163.7661506121285
$ws.Range("B11").Value = 1672.0183886859197
$ws.Range("C11").Value = 1602.2163373852031
$ws.Range("A12").Value = 2859.3300321679417
$ws.Range("B12").Value = 2220.3401154197459
$ws.Range("C12").Value = 1930.1062011589634
$ws.Range("A13").Value = 2573.4500117771727
$ws.Range("B13").Value = 2069.4587285737271
$ws.Range("C13").Value = 2070.9582771484247
$ws.Range("A14").Value = 2609.3935857794863
$ws.Range("B14").Value = 2264.4745278119372
$ws.Range("C14").Value = 2131.8727646395555
$ws.Range("A15").Value = 2574.253237504885
$ws.Range("B15").Value = 1943.0136835459234
$ws.Range("C15").Value = 2182.3893458170587
$ws.Range("A16").Value = 2591.9527981864389
$ws.Range("B16").Value = 1796.809889737392
$ws.Range("C16").Value = 1568.4260851246033

$wb.Save()
